# Fix i/o errors and discard MIAPPE-based templates
#
# Rename the "Source Name" input/output columns of the
# Events-CropResidueIncorporation annotation table to "Sample Name",
# matching the ISA "Sample Name" vocabulary used elsewhere in the
# FAIRagro templates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events-CropResidueIncorporation")

# Update the header cells directly …
$ws.Range("A1").Value = "Input [Sample Name]"
$ws.Range("AK1").Value = "Output [Sample Name]"

# … and keep the ListObject / table column names (which drive
# xl/tables/table1.xml) in sync with the header text.
$table = $ws.ListObjects.Item("annotationTable")
$table.ListColumns.Item(1).Name = "Input [Sample Name]"
$table.ListColumns.Item(37).Name = "Output [Sample Name]"
